$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Role class swap: columns A4/A5 now hold Rouge/Ranger (previously Ranger/Rouge) ---
$ws.Range("A4").Value = "Rouge"
$ws.Range("A5").Value = "Ranger"

# --- New column headers (row 1) ---
# Existing: A1 = "name"
# (order matters for shared-string table index assignment, so the
#  "increment" headers are entered before the "starting max" headers)
$ws.Range("F1").Value = "hp increment"
$ws.Range("G1").Value = "dmg increment"
$ws.Range("H1").Value = "def increment"
$ws.Range("I1").Value = "stamina increment"
$ws.Range("B1").Value = "starting max hp"
$ws.Range("C1").Value = "starting max dmg"
$ws.Range("D1").Value = "starting def"
$ws.Range("E1").Value = "starting max stamina"

# --- New per-class stat columns (B:I) for rows 2-7 ---
# order: starting hp, starting dmg, starting def, starting stamina,
#        hp inc, dmg inc, def inc, stamina inc
$classData = @{
    "Warrior" = @(12, 4, 2, 10, 2, 1, 1, 3)
    "Mage"    = @(6, 12, 0, 5, 1, 3, 1, 1)
    "Rouge"   = @(9, 7, 1, 7, 2, 1, 2, 2)
    "Ranger"  = @(8, 8, 1, 6, 1, 2, 2, 1)
    "Acolyte" = @(6, 6, 1, 10, 1, 2, 1, 3)
    "Slayer"  = @(8, 10, 0, 10, 1, 3, 0, 1)
}

for ($row = 2; $row -le 7; $row++) {
    $name = $ws.Cells.Item($row, 1).Value()
    $stats = $classData[$name]
    for ($i = 0; $i -lt $stats.Length; $i++) {
        $ws.Cells.Item($row, $i + 2).Value = $stats[$i]
    }
}

# --- Autofit the new columns to match the resulting content widths ---
$null = $ws.Columns("B:I").AutoFit()

# --- Update the active selection to match the saved view ---
$null = $ws.Range("F12").Select()
